# 2021/03/30 16:58 — new results.
# Update the two revised expert scores (rows 24 and 29, column B) and
# reset the sheet's zoom level back to 100%. The AVERAGE formula in B32
# recalculates automatically from the updated inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24").Value = 0.43209999999999998
$ws.Range("B29").Value = 0.38929999999999998

$excel.ActiveWindow.Zoom = 100
